$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 85: Argent (silver) price update, mirrors the prior row's values
# except for the date. All columns are stored as text in this sheet, so
# force text formatting before assigning the literal strings to avoid
# Excel auto-coercing the date-like / numeric-like text into real
# dates/numbers.
$row = $ws.Range("A85:J85")
$row.NumberFormat = "@"

$ws.Range("A85").Value = "2025-05-25"
$ws.Range("B85").Value = "35.5"
$ws.Range("C85").Value = "35.4"
$ws.Range("D85").Value = "0.94"
$ws.Range("E85").Value = "0.258"
$ws.Range("F85").Value = "0.09"
$ws.Range("G85").Value = "5,373"
$ws.Range("H85").Value = "8,045"
$ws.Range("I85").Value = "8,095"
$ws.Range("J85").Value = "7.2241"

# Restore the default cell style so the new row doesn't carry an explicit
# number-format style index (matches the plain, style-less cells used by
# every other data row in the sheet).
$row.Style = "Normal"
